$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Remove the obsolete screens (rows) from the translation table.
# Rows 38-52 are deleted entirely, which shifts the former rows 53-56
# up to become the new rows 38-41.
$ws.Range("A38:A52").EntireRow.Delete()
